# Weekly update: a new week's price record is inserted at row 19,
# pushing the existing rows 19-109 down to 20-110.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19 (shifts rows 19..109 down to 20..110,
# carrying their formatting - including the date-formatted column D - with them).
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly record.
$ws.Cells.Item(19, 1).Value = 8
$ws.Cells.Item(19, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(19, 3).Value = "Coquimbo"
$ws.Cells.Item(19, 4).Value = 44561
$ws.Cells.Item(19, 5).Value = 4
$ws.Cells.Item(19, 6).Value = 100112001
$ws.Cells.Item(19, 7).Value = "Berenjena"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 520
$ws.Cells.Item(19, 11).Value = 9000
$ws.Cells.Item(19, 12).Value = 10000
$ws.Cells.Item(19, 13).Value = 9500
$ws.Cells.Item(19, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(19, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(19, 16).Value = 158
$ws.Cells.Item(19, 17).Value = 60
$ws.Cells.Item(19, 18).Value = "Hortaliza"
